# Add two new worksheets, "contingencia_abs" and "contingencia_rel",
# after the existing "frec_ventas" sheet, each holding a CIUDAD x
# (MEDIANA / MICROEMPRESA / PEQUEÑA) contingency table.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsAbs = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsAbs.Name = "contingencia_abs"

$wsRel = $wb.Worksheets.Add([Type]::Missing, $wsAbs)
$wsRel.Name = "contingencia_rel"

$header = @("CIUDAD", "MEDIANA        ", "MICROEMPRESA   ", "PEQUEÑA        ")

$cities = @(
    "DAULE                                             ",
    "ELOY ALFARO                                       ",
    "GUAYAQUIL                                         ",
    "MILAGRO                                           ",
    "NARANJITO                                         ",
    "SAMBORONDÓN                                       ",
    "SANTA LUCIA                                       ",
    "VELASCO IBARRA                                    "
)

$absRows = @(
    @(1, 1, 0),
    @(0, 2, 0),
    @(2, 117, 28),
    @(0, 2, 0),
    @(0, 4, 0),
    @(0, 0, 1),
    @(0, 1, 0),
    @(0, 1, 1)
)

$relRows = @(
    @(0.62, 0.62, 0.0),
    @(0.0, 1.24, 0.0),
    @(1.24, 72.67, 17.39),
    @(0.0, 1.24, 0.0),
    @(0.0, 2.48, 0.0),
    @(0.0, 0.0, 0.62),
    @(0.0, 0.62, 0.0),
    @(0.0, 0.62, 0.62)
)

for ($c = 0; $c -lt 4; $c++) {
    $wsAbs.Cells.Item(1, $c + 1).Value = $header[$c]
    $wsRel.Cells.Item(1, $c + 1).Value = $header[$c]
}

for ($r = 0; $r -lt $cities.Count; $r++) {
    $rowNum = $r + 2

    $wsAbs.Cells.Item($rowNum, 1).Value = $cities[$r]
    $wsRel.Cells.Item($rowNum, 1).Value = $cities[$r]

    for ($c = 0; $c -lt 3; $c++) {
        $wsAbs.Cells.Item($rowNum, $c + 2).Value = $absRows[$r][$c]
        $wsRel.Cells.Item($rowNum, $c + 2).Value = $relRows[$r][$c]
    }
}
